$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Writes $Text into $Range as literal text. Excel's COM layer auto-detects
# plain numeric-looking strings (e.g. "187.77", "1.00") and silently stores
# them as numbers, which would corrupt values like "1.00" (-> 1) or
# "12.30" (-> 12.3). Prefixing with an apostrophe forces text, matching
# what a user typing '187.77 into a cell gets; the cell's style is then put
# back to the default so only the value (not formatting) changes.
function Set-TextValue {
    param($Range, [string]$Text)
    if ($Text -match '^[+-]?[0-9]+(\.[0-9]+)*$') {
        $Range.Value = "'" + $Text
        $Range.Style = 'Normal'
    } else {
        $Range.Value = $Text
    }
}

Set-TextValue $ws.Range('D2') '75.016.27'
Set-TextValue $ws.Range('E2') '  +1.79%  '
Set-TextValue $ws.Range('D3') '2.821.82'
Set-TextValue $ws.Range('E3') '  +7.73%  '
Set-TextValue $ws.Range('E4') '  +0.06%  '
Set-TextValue $ws.Range('D5') '187.77'
Set-TextValue $ws.Range('E5') '  +2.24%  '
Set-TextValue $ws.Range('D6') '594.34'
Set-TextValue $ws.Range('E6') '  +2.29%  '
Set-TextValue $ws.Range('E7') '  +0.02%  '
Set-TextValue $ws.Range('D8') '0.553'
Set-TextValue $ws.Range('E8') '  +3.93%  '
Set-TextValue $ws.Range('D9') '0.192'
Set-TextValue $ws.Range('E9') '  -2.24%  '
Set-TextValue $ws.Range('D10') '2.820.81'
Set-TextValue $ws.Range('E10') '  +7.67%  '
Set-TextValue $ws.Range('E11') '  -0.97%  '
Set-TextValue $ws.Range('D12') '0.371'
Set-TextValue $ws.Range('E12') '  +3.78%  '
Set-TextValue $ws.Range('D13') '4.84'
Set-TextValue $ws.Range('E13') '  +2.31%  '
Set-TextValue $ws.Range('D14') '3.339.05'
Set-TextValue $ws.Range('E14') '  +8.54%  '
Set-TextValue $ws.Range('D15') '75.012.96'
Set-TextValue $ws.Range('E15') '  +1.96%  '
Set-TextValue $ws.Range('D16') '0.0000187'
Set-TextValue $ws.Range('E16') '  -0.02%  '
Set-TextValue $ws.Range('D17') '26.86'
Set-TextValue $ws.Range('E17') '  +3.45%  '
Set-TextValue $ws.Range('D18') '2.821.54'
Set-TextValue $ws.Range('E18') '  +7.73%  '
Set-TextValue $ws.Range('D19') '8.88'
Set-TextValue $ws.Range('E19') '  -0.72%  '
Set-TextValue $ws.Range('D20') '12.30'
Set-TextValue $ws.Range('E20') '  +4.28%  '
Set-TextValue $ws.Range('D21') '375.96'
Set-TextValue $ws.Range('E21') '  +1.45%  '
Set-TextValue $ws.Range('D22') '2.24'
Set-TextValue $ws.Range('E22') '  +0.46%  '
Set-TextValue $ws.Range('D23') '4.07'
Set-TextValue $ws.Range('E23') '  +0.23%  '
Set-TextValue $ws.Range('D24') '6.15'
Set-TextValue $ws.Range('E24') '  -0.89%  '
Set-TextValue $ws.Range('D26') '71.01'
Set-TextValue $ws.Range('E26') '  +2.34%  '
Set-TextValue $ws.Range('D28') '4.16'
Set-TextValue $ws.Range('E28') '  +1.40%  '
Set-TextValue $ws.Range('D29') '9.58'
Set-TextValue $ws.Range('E29') '  +4.06%  '
Set-TextValue $ws.Range('D30') '0.0000103'
Set-TextValue $ws.Range('E30') '  +11.20%  '
Set-TextValue $ws.Range('D31') '1.00'
Set-TextValue $ws.Range('E31') '  -0.15%  '
Set-TextValue $ws.Range('D32') '511.34'
Set-TextValue $ws.Range('E32') '  -0.77%  '
Set-TextValue $ws.Range('D33') '1.38'
Set-TextValue $ws.Range('E33') '  +0.72%  '
Set-TextValue $ws.Range('D34') '7.73'
Set-TextValue $ws.Range('E34') '  +1.76%  '
Set-TextValue $ws.Range('E35') '  +3.90%  '
Set-TextValue $ws.Range('E36') '  +0.02%  '
Set-TextValue $ws.Range('D37') '162.69'
Set-TextValue $ws.Range('E37') '  +0.73%  '
Set-TextValue $ws.Range('D38') '20.02'
Set-TextValue $ws.Range('E38') '  +4.90%  '
Set-TextValue $ws.Range('E39') '  -0.86%  '
Set-TextValue $ws.Range('D40') '19.42'
Set-TextValue $ws.Range('E40') '  +0.81%  '
Set-TextValue $ws.Range('D41') '184.68'
Set-TextValue $ws.Range('E41') '  +17.36%  '
Set-TextValue $ws.Range('B43') 'PolygonEcosystemToken'
Set-TextValue $ws.Range('C43') 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue $ws.Range('D43') '0.342'
Set-TextValue $ws.Range('E43') '  +5.47%  '
Set-TextValue $ws.Range('B44') 'RenderToken'
Set-TextValue $ws.Range('C44') 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-TextValue $ws.Range('D44') '5.03'
Set-TextValue $ws.Range('E44') '  +3.97%  '
Set-TextValue $ws.Range('D45') '1.67'
Set-TextValue $ws.Range('E45') '  +0.43%  '
Set-TextValue $ws.Range('D46') '1.21'
Set-TextValue $ws.Range('E46') '  +2.87%  '
Set-TextValue $ws.Range('D47') '40.00'
Set-TextValue $ws.Range('E47') '  +3.68%  '
Set-TextValue $ws.Range('E48') '  +0.19%  '
Set-TextValue $ws.Range('D49') '0.0855'
Set-TextValue $ws.Range('E49') '  -1.97%  '
Set-TextValue $ws.Range('E50') '  +8.94%  '
Set-TextValue $ws.Range('D51') '3.72'
Set-TextValue $ws.Range('E51') '  +3.37%  '
